$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cadastro")
$ws2 = $wb.Worksheets.Item("Pesquisa")

# --- Sheet "Pesquisa" (Search) updates ---
# New product entry: Bose SoundLink Wireless Speaker / Speakers
$ws2.Range("B2").Value = "Bose SoundLink Wireless Speaker"
$ws2.Range("A2").Value = "Speakers"

# Resize column B to fit the new, longer product name
$ws2.Range("B:B").ColumnWidth = 29.83

# Second product entry further down the sheet
$ws2.Range("I15").Value = "HP USB 3 Button Optical Mouse"

# Bring over the same "empty but formatted" cell that exists at H12,
# landing on B4 (mirrors the underline-styled blank cell from Cadastro!A3)
$ws1.Range("A3").Copy() | Out-Null
$ws2.Range("B4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Sheet "Cadastro" (Registration) updates ---
# New "Pais"/"Brazil" column
$ws1.Range("K1").Value = "Pais"
$ws1.Range("K2").Value = "Brazil"

# Bring over the same "empty but formatted" cell that exists at H12 (Pesquisa),
# landing on K3
$ws2.Range("H12").Copy() | Out-Null
$ws1.Range("K3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Username corrected/renamed
$ws1.Range("A2").Value = "abdielCordeiro"

# --- Selection / active cell bookkeeping ---
$ws2.Range("I15").Select() | Out-Null
$ws1.Select() | Out-Null
$ws1.Range("A2").Select() | Out-Null
